# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the fe03d05e-fd73-4c39-a5c4-bac138d75e85 file across the Overview,
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-17 04:43:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-17 04:43:22"
$zhcn.Range("K3").Value = "2016-08-17 04:43:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-17 04:43:27"
$dede.Range("K3").Value = "2016-08-17 04:43:47"
